# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 211.19
$wsSummary.Range("E3").Value = 114.2

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 4
$wsRepay.Range("F4").Value = 921.65
$wsRepay.Range("G4").Value = 3211.57

# Row 5 - H5 used to be a formula (=G4*(12%/365)*B5) styled with a dedicated
# number-format; it becomes a plain value using the same style as the rest
# of the column (copy formats from H4 first, then overwrite the value).
$wsRepay.Range("F5").Value = 932.09
$wsRepay.Range("G5").Value = 2279.48
$wsRepay.Range("H4").Copy()
$wsRepay.Range("H5").PasteSpecial(-4122)
$wsRepay.Range("H5").Value = 31.68

# Row 6
$wsRepay.Range("F6").Value = 940.54
$wsRepay.Range("G6").Value = 1338.94
$wsRepay.Range("H6").Value = 23.23

# Row 7
$wsRepay.Range("F7").Value = 950.56
$wsRepay.Range("G7").Value = 388.38
$wsRepay.Range("H7").Value = 13.21

# Row 8 - D8/E8 did not previously exist; add them (formatted like D7:E7,
# left empty) then update the remaining values.
$wsRepay.Range("D7:E7").Copy()
$wsRepay.Range("D8:E8").PasteSpecial(-4122)
$wsRepay.Range("F8").Value = 388.38
$wsRepay.Range("H8").Value = 3.96
$wsRepay.Range("K8").Value = 392.34
$wsRepay.Range("P8").Value = 392.34

# Column O (rows 2-8) is dropped entirely.
$wsRepay.Range("O2").Clear()
$wsRepay.Range("O3").Clear()
$wsRepay.Range("O4").Clear()
$wsRepay.Range("O5").Clear()
$wsRepay.Range("O6").Clear()
$wsRepay.Range("O7").Clear()
$wsRepay.Range("O8").Clear()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

$wsTrans.Range("A2").Value = 59
$wsTrans.Range("A3").Value = 57

# ---------------------------------------------------------------------------
# Selections - applied last, in tab order, so the final selection/active
# sheet state mirrors the saved workbook (Transactions stays the active tab).
# ---------------------------------------------------------------------------
$wsSummary.Range("C4").Select()
$wsRepay.Range("A9:XFD9").Select()
$wsTrans.Range("A2:L3").Select()
